$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13, shifting rows 13-24 down to 14-25.
$ws.Rows("13:13").Insert()

# The insert duplicated column-A formatting into the new row 13; row 13 in the
# target layout has no entry in column A at all, so clear it back out.
$ws.Range("A13").Clear()

# Row 13 (new): Docentes responsaveis value moves up here (B/C only).
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = '1097178 - João Batista de Almeida e Silva'
$ws.Range("C13").Value = '1097178 - João Batista de Almeida e Silva'

# Row 10: Objetivos text is replaced with the new Portuguese objectives paragraph.
$ws.Range("B10").Value = 'Proporcionar ao estudante conhecimentos práticos nos processos tecnológicos de preparação de bebidas fermentadas e destiladas.'
$ws.Range("C10").Value = 'Proporcionar ao estudante conhecimentos práticos nos processos tecnológicos de preparação de bebidas fermentadas e destiladas.'

# Row 14: Programa resumido text.
$ws.Range("B14").Value = 'Elaboração prática de cerveja, cachaça, fermentados e destilados de frutas, cereais e tuberculos, vinhos e análise sensorial.'
$ws.Range("C14").Value = 'Elaboração prática de cerveja, cachaça, fermentados e destilados de frutas, cereais e tuberculos, vinhos e análise sensorial.'

# Row 16: Programa (long Portuguese syllabus) text.
$ws.Range("B16").Value = '1. Elaboração de cerveja: matérias-primas, preparação do mosto, tecnologia de fermentação e maturação.2. Elaboração de aguardente: matérias-primas, preparação do mosto, tecnologia de fermentação, destilação, maturação.3. Elaboração de destilados de frutas: matérias-primas, preparação do mosto, tecnologia de fermentação, destilação, maturação.4. Elaboração e vinhos: matérias-primas, preparação do mosto, tecnologia de fermentação, maturação.5. Análise sensorial: teste sensorial das bebidas preparadas nos itens anteriores'
$ws.Range("C16").Value = '1. Elaboração de cerveja: matérias-primas, preparação do mosto, tecnologia de fermentação e maturação.2. Elaboração de aguardente: matérias-primas, preparação do mosto, tecnologia de fermentação, destilação, maturação.3. Elaboração de destilados de frutas: matérias-primas, preparação do mosto, tecnologia de fermentação, destilação, maturação.4. Elaboração e vinhos: matérias-primas, preparação do mosto, tecnologia de fermentação, maturação.5. Análise sensorial: teste sensorial das bebidas preparadas nos itens anteriores'

# Row 19: Metodo text.
$ws.Range("B19").Value = 'Relatórios e seminários sobre os experimentos'
$ws.Range("C19").Value = 'Relatórios e seminários sobre os experimentos'

# Row 20: Criterio text.
$ws.Range("B20").Value = 'Média aritmética entre os relatórios e seminários'
$ws.Range("C20").Value = 'Média aritmética entre os relatórios e seminários'

# Row 21: Norma de recuperacao text.
$ws.Range("B21").Value = 'A recuperação será feita por meio de prova escrita (PR) e a média final (MF) será calculada pela equação: MF = (NF + PR)/2.'
$ws.Range("C21").Value = 'A recuperação será feita por meio de prova escrita (PR) e a média final (MF) será calculada pela equação: MF = (NF + PR)/2.'

# Row 22: Bibliografia text.
$ws.Range("B22").Value = '1. AQUARONE, E.; BORZANI, W.; SCHMIDELL, W.; LIMA, U. A. Biotecnologia na Produção deAlimentos. V. 4, Biotecnologia Industrial, São Paulo: Edgard Blücher Ltda. 2001.4. DUVAL, G. Fabricação de Vinhos de Frutas. S.I.A.RJ:Ministério da Agricultura, 1947.5. HOUGH, J.S. Biotecnología de La cerveza y de la malta. Editorial ACRIBA S/A, 1978.6. LIMA, U. A. Aguardente: fabricação em pequenas destilarias. Ed. FEALQ. 1999.7. MARTINELLI FILHO, A. Tecnologia de Vinhos e Vinagres de Frutas. Agroindústria de BaixoInvestimento. Departamento de Tecnologia Rural da ESALQ/USP.8. MORRETO, E. et al. Vinhos e Vinagres: Processamento e Análises. FlorianópolisEditoraUFSC, 1988.9. PACHECO, A. O. Manual do Bar. São Paulo. Editora SENAC, 1996.10. STANIER, R. Y.; INGRAHAM, J. L., WHEELIS, M. L.; PAINTER, P. R. The Microbial World.Englewood Cliffs, New Jersey, 1986.11.Venturini Filho, W.G. Bebidas Alcoólicas. Ciência e Tecnologia. São Paulo. Edgar Blucher Ltda. 2a. Edição. 2016. 575 p.'
$ws.Range("C22").Value = '1. AQUARONE, E.; BORZANI, W.; SCHMIDELL, W.; LIMA, U. A. Biotecnologia na Produção deAlimentos. V. 4, Biotecnologia Industrial, São Paulo: Edgard Blücher Ltda. 2001.4. DUVAL, G. Fabricação de Vinhos de Frutas. S.I.A.RJ:Ministério da Agricultura, 1947.5. HOUGH, J.S. Biotecnología de La cerveza y de la malta. Editorial ACRIBA S/A, 1978.6. LIMA, U. A. Aguardente: fabricação em pequenas destilarias. Ed. FEALQ. 1999.7. MARTINELLI FILHO, A. Tecnologia de Vinhos e Vinagres de Frutas. Agroindústria de BaixoInvestimento. Departamento de Tecnologia Rural da ESALQ/USP.8. MORRETO, E. et al. Vinhos e Vinagres: Processamento e Análises. FlorianópolisEditoraUFSC, 1988.9. PACHECO, A. O. Manual do Bar. São Paulo. Editora SENAC, 1996.10. STANIER, R. Y.; INGRAHAM, J. L., WHEELIS, M. L.; PAINTER, P. R. The Microbial World.Englewood Cliffs, New Jersey, 1986.11.Venturini Filho, W.G. Bebidas Alcoólicas. Ciência e Tecnologia. São Paulo. Edgar Blucher Ltda. 2a. Edição. 2016. 575 p.'
